$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to remain text so numeric-looking values are not
# silently converted to Excel numbers (matches original inline-string format).
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "63.874.18"
$ws.Range("E2").Value = "  +0.00%  "

# Row 3
$ws.Range("D3").Value = "2.629.95"
$ws.Range("E3").Value = "  +0.21%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "595.79"
$ws.Range("E5").Value = "  -0.05%  "

# Row 6
$ws.Range("D6").Value = "151.76"
$ws.Range("E6").Value = "  +1.12%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
$ws.Range("D9").Value = "0.114"
$ws.Range("E9").Value = "  +4.13%  "

# Row 10
$ws.Range("D10").Value = "5.85"
$ws.Range("E10").Value = "  +3.64%  "

# Row 11
$ws.Range("D11").Value = "0.398"
$ws.Range("E11").Value = "  +4.29%  "

# Row 12
$ws.Range("E12").Value = "  +0.98%  "

# Row 13
$ws.Range("D13").Value = "28.17"
$ws.Range("E13").Value = "  +1.80%  "

# Row 14
$ws.Range("D14").Value = "3.103.30"
$ws.Range("E14").Value = "  +0.26%  "

# Row 15
$ws.Range("D15").Value = "63.728.96"

# Row 16
$ws.Range("E16").Value = "  +10.49%  "

# Row 17
$ws.Range("D17").Value = "2.647.58"
$ws.Range("E17").Value = "  +1.85%  "

# Row 18
$ws.Range("D18").Value = "12.20"
$ws.Range("E18").Value = "  -0.24%  "

# Row 19
$ws.Range("D19").Value = "4.81"
$ws.Range("E19").Value = "  +4.27%  "

# Row 20
$ws.Range("D20").Value = "348.31"
$ws.Range("E20").Value = "  -0.32%  "

# Row 21
$ws.Range("D21").Value = "7.03"
$ws.Range("E21").Value = "  +1.89%  "

# Row 22
$ws.Range("E22").Value = "  +0.24%  "

# Row 23
$ws.Range("D23").Value = "67.58"
$ws.Range("E23").Value = "  +1.97%  "

# Row 24
$ws.Range("E24").Value = "  -2.23%  "

# Row 25
$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D25").Value = "1.68"
$ws.Range("E25").Value = "  +0.40%  "

# Row 26
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "9.17"
$ws.Range("E26").Value = "  -0.34%  "

# Row 27
$ws.Range("D27").Value = "8.40"
$ws.Range("E27").Value = "  +2.17%  "

# Row 28
$ws.Range("D28").Value = "555.41"
$ws.Range("E28").Value = "  +1.38%  "

# Row 29
$ws.Range("E29").Value = "  -1.51%  "

# Row 30
$ws.Range("E30").Value = "  +0.03%  "

# Row 31
$ws.Range("E31").Value = "  +1.26%  "

# Row 32
$ws.Range("D32").Value = "0.0₃0896"
$ws.Range("E32").Value = "  +5.86%  "

# Row 33
$ws.Range("D33").Value = "1.82"
$ws.Range("E33").Value = "  +5.02%  "

# Row 34
$ws.Range("D34").Value = "5.40"
$ws.Range("E34").Value = "  +3.40%  "

# Row 35
$ws.Range("D35").Value = "6.14"
$ws.Range("E35").Value = "  +0.87%  "

# Row 36
$ws.Range("E36").Value = "  -2.56%  "

# Row 37
$ws.Range("E37").Value = "  +2.67%  "

# Row 38
$ws.Range("D38").Value = "2.00"
$ws.Range("E38").Value = "  +2.29%  "

# Row 39
$ws.Range("D39").Value = "19.96"
$ws.Range("E39").Value = "  +3.09%  "

# Row 40
$ws.Range("E40").Value = "  -0.07%  "

# Row 41
$ws.Range("E41").Value = "  -0.01%  "

# Row 42
$ws.Range("D42").Value = "168.44"
$ws.Range("E42").Value = "  -1.16%  "

# Row 43
$ws.Range("E43").Value = "  +3.65%  "

# Row 44
$ws.Range("D44").Value = "23.59"
$ws.Range("E44").Value = "  +9.93%  "

# Row 45
$ws.Range("D45").Value = "0.0587"
$ws.Range("E45").Value = "  -1.51%  "

# Row 46
$ws.Range("E46").Value = "  +10.21%  "

# Row 47
$ws.Range("D47").Value = "0.639"
$ws.Range("E47").Value = "  +1.63%  "

# Row 48
$ws.Range("E48").Value = "  +2.93%  "

# Row 49
$ws.Range("D49").Value = "0.0972"
$ws.Range("E49").Value = "  +0.68%  "

# Row 50
$ws.Range("D50").Value = "19.34"
$ws.Range("E50").Value = "  +1.08%  "

# Row 51
$ws.Range("D51").Value = "0.0₆0234"
$ws.Range("E51").Value = "  +19.56%  "
